$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.713.02"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.460.93"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'571.19"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'147.07"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'5.17"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").Value = "'28.85"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "2.908.74"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "62.585.39"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "2.467.28"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'7.64"
$ws.Range("E18").Value = "  -6.36%  "
$ws.Range("D19").Value = "'10.72"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("D20").Value = "'2.32"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'321.08"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'10.15"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "'64.75"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("D26").Value = "'640.48"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -5.30%  "
$ws.Range("D31").Value = "'7.88"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "'150.27"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "0.0₆0306"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "'153.54"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").Value = "'15.40"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.604"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'20.16"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "'0.0901"
$ws.Range("E51").Value = "  -1.90%  "
